# Updates the "Price" (column D) and "Volume(1h)" (column E) columns of the
# cryptocurrency listing on the active worksheet to the refreshed values
# captured by the latest GitHub Actions run.
#
# All of these cells hold plain *text* in the original workbook (t="inlineStr"),
# even though many of the new values look like plain numbers (e.g. "1.000",
# "43.84"). Assigning such a string straight to Range.Value causes Excel to
# auto-coerce it into a real number (dropping the trailing zero, turning
# "328.52" into 328.52, etc.), which would corrupt the displayed text. To
# avoid that, for any replacement value that parses as a number we briefly
# force the cell's number format to Text ("@") before writing the value, then
# clear the format back to the sheet's default styling afterwards so no
# visual/style changes are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = [ordered]@{
    'D2' = '27.955.65'
    'E2' = '  +0.89%  '
    'D3' = '1.763.72'
    'E3' = '  -0.73%  '
    'D4' = '1.000'
    'E4' = '  -0.06%  '
    'D5' = '328.52'
    'E5' = '  +0.64%  '
    'E6' = '  -0.06%  '
    'D7' = '0.4645'
    'E7' = '  +0.64%  '
    'D8' = '0.3517'
    'E8' = '  -2.09%  '
    'D9' = '43.84'
    'E9' = '  +4.45%  '
    'D10' = '0.07361'
    'E10' = '  -1.58%  '
    'D11' = '1.085'
    'E11' = '  -1.55%  '
    'D12' = '0.9994'
    'E12' = '  -0.12%  '
    'D13' = '20.61'
    'E13' = '  -1.15%  '
    'D14' = '5.997'
    'E14' = '  -0.65%  '
    'D15' = '7.158'
    'E15' = '  -1.05%  '
    'D16' = '1.763.03'
    'E16' = '  -0.69%  '
    'D17' = '92.50'
    'E17' = '  -1.18%  '
    'D18' = '0.00001053'
    'E18' = '  -0.49%  '
    'D19' = '0.06423'
    'E19' = '  +0.13%  '
    'D20' = '0.9996'
    'D21' = '16.84'
    'E21' = '  -1.44%  '
    'D22' = '5.766'
    'E22' = '  -0.35%  '
    'D23' = '27.964.44'
    'E23' = '  +0.65%  '
    'D24' = '11.15'
    'E24' = '  -1.27%  '
    'D25' = '2.153'
    'E25' = '  +3.55%  '
    'D26' = '162.70'
    'E26' = '  -1.09%  '
    'D27' = '20.02'
    'E27' = '  -1.71%  '
    'D28' = '1.966.18'
    'E28' = '  -0.68%  '
    'D29' = '2.169'
    'E29' = '  +0.20%  '
    'D30' = '123.00'
    'E30' = '  -2.55%  '
    'D31' = '1.067'
    'E31' = '  -2.74%  '
    'D32' = '0.09280'
    'E32' = '  +0.57%  '
    'D33' = '3.644'
    'E33' = '  -0.89%  '
    'D34' = '5.549'
    'E34' = '  +0.20%  '
    'D35' = '11.67'
    'E35' = '  -1.13%  '
    'D36' = '0.02271'
    'E36' = '  -1.03%  '
    'E37' = '  -0.65%  '
    'D38' = '0.2061'
    'E38' = '  -1.33%  '
    'D39' = '4.908'
    'E39' = '  -1.09%  '
    'D40' = '0.6137'
    'E40' = '  -2.77%  '
    'D41' = '1.180'
    'E41' = '  -0.07%  '
    'D42' = '7.794'
    'E42' = '  +0.23%  '
    'D43' = '1.360'
    'E43' = '  -2.42%  '
    'D44' = '13.12'
    'E44' = '  -1.10%  '
    'D45' = '3.736'
    'E45' = '  +0.17%  '
    'D46' = '0.5788'
    'E46' = '  -1.72%  '
    'D47' = '122.77'
    'E47' = '  +0.38%  '
    'D48' = '1.927'
    'E48' = '  -1.20%  '
    'D49' = '0.06814'
    'E49' = '  -1.86%  '
    'D50' = '1.123'
    'E50' = '  -1.40%  '
    'D51' = '72.12'
    'E51' = '  -0.22%  '
}

foreach ($addr in $values.Keys) {
    $newValue = $values[$addr]
    $cell = $ws.Range($addr)

    # A simple "does this look like a plain number" check (optional
    # sign, digits, optional single decimal point). Values containing a
    # second '.' (e.g. "27.955.65") or a '%'/spaces (e.g. "  +0.89%  ")
    # will not match and are therefore safe to assign as-is.
    $looksNumeric = $newValue -match '^\s*[-+]?\d+(\.\d+)?\s*$'

    if ($looksNumeric) {
        # Force text storage so Excel keeps the exact string (e.g. "1.000",
        # "328.52") instead of silently converting it to a number.
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        # Restore default cell styling/format now that the text is locked in.
        $cell.ClearFormats()
    } else {
        # Values such as "  +0.89%  " or "27.955.65" (two dots) can never be
        # parsed as a number by Excel, so they are safe to assign directly.
        $cell.Value = $newValue
    }
}
